$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3477.2727
$ws.Range("J17").Value = 1218.3907
$ws.Range("L17").Value = 3655.1721
$ws.Range("N17").Value = -3991.1721
$ws.Range("H70").Value = 2335.1853
$ws.Range("I70").Value = 993.75
$ws.Range("J70").Value = 2900
$ws.Range("K70").Value = 2981.25
$ws.Range("L70").Value = 8700
$ws.Range("M70").Value = -2711.25
$ws.Range("N70").Value = -9240
$ws.Range("H73").Value = 2335.1853
$ws.Range("I73").Value = 993.75
$ws.Range("J73").Value = 2900
$ws.Range("K73").Value = 2981.25
$ws.Range("L73").Value = 8700
$ws.Range("M73").Value = -2045.25
$ws.Range("N73").Value = -10572
$ws.Range("H94").Value = 2963.6365
$ws.Range("I94").Value = 2840
$ws.Range("K94").Value = 2840
$ws.Range("M94").Value = -2389
$ws.Range("H100").Value = 4754.522
$ws.Range("I100").Value = 4160.778
$ws.Range("J100").Value = 6892
$ws.Range("K100").Value = 4160.778
$ws.Range("L100").Value = 6892
$ws.Range("M100").Value = -3619.778
$ws.Range("N100").Value = -7974
$ws.Range("H107").Value = 536.88464
$ws.Range("I107").Value = 618.6111
$ws.Range("J107").Value = 353
$ws.Range("K107").Value = 618.6111
$ws.Range("L107").Value = 353
$ws.Range("M107").Value = 1301.3889
$ws.Range("N107").Value = -4193
$ws.Range("H116").Value = 4666.143
$ws.Range("I116").Value = 3700
$ws.Range("J116").Value = 5632.2856
$ws.Range("K116").Value = 3700
$ws.Range("L116").Value = 5632.2856
$ws.Range("M116").Value = -258
$ws.Range("N116").Value = -12516.2856
$ws.Range("H137").Value = 4549850
$ws.Range("I137").Value = 5886847.5
$ws.Range("K137").Value = 17660542.5
$ws.Range("M137").Value = -17657992.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5774.488
$ws.Range("I32").Value = 4187.418
$ws.Range("K32").Value = 4187.418
$ws.Range("M32").Value = -3900.418
$ws.Range("H45").Value = 1537.3112
$ws.Range("I45").Value = 1066.5714
$ws.Range("J45").Value = 3184.9
$ws.Range("K45").Value = 1066.5714
$ws.Range("L45").Value = 3184.9
$ws.Range("M45").Value = -689.5714
$ws.Range("N45").Value = -3938.9
$ws.Range("H74").Value = 575.88
$ws.Range("I74").Value = 544.6818
$ws.Range("K74").Value = 544.6818
$ws.Range("M74").Value = 329.3182
$ws.Range("H77").Value = 575.88
$ws.Range("I77").Value = 544.6818
$ws.Range("K77").Value = 2723.409
$ws.Range("M77").Value = 1644.591
$ws.Range("H102").Value = 4412.857
$ws.Range("I102").Value = 3678
$ws.Range("K102").Value = 3678
$ws.Range("M102").Value = -2056
$ws.Range("H110").Value = 1293.1765
$ws.Range("I110").Value = 640.7143
$ws.Range("J110").Value = 2347.1538
$ws.Range("K110").Value = 640.7143
$ws.Range("L110").Value = 2347.1538
$ws.Range("M110").Value = 1404.2857
$ws.Range("N110").Value = -6437.1538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3007.7778
$ws.Range("I99").Value = 1861.6666
$ws.Range("J99").Value = 5300
$ws.Range("K99").Value = 1861.6666
$ws.Range("L99").Value = 5300
$ws.Range("M99").Value = -363.6666
$ws.Range("N99").Value = -8296
$ws.Range("H105").Value = 1825.5807
$ws.Range("I105").Value = 1490
$ws.Range("J105").Value = 2067.9443
$ws.Range("K105").Value = 1490
$ws.Range("L105").Value = 2067.9443
$ws.Range("M105").Value = 257
$ws.Range("N105").Value = -5561.9443

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2276736.8
$ws.Range("I31").Value = 2705579
$ws.Range("J31").Value = 10000
$ws.Range("K31").Value = 2705579
$ws.Range("L31").Value = 10000
$ws.Range("M31").Value = -2705284
$ws.Range("N31").Value = -10590
$ws.Range("H34").Value = 2276736.8
$ws.Range("I34").Value = 2705579
$ws.Range("J34").Value = 10000
$ws.Range("K34").Value = 2705579
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = -2705377
$ws.Range("N34").Value = -10404
$ws.Range("H58").Value = 8476651
$ws.Range("I58").Value = 1436.1842
$ws.Range("J58").Value = 23812754
$ws.Range("K58").Value = 1436.1842
$ws.Range("L58").Value = 23812754
$ws.Range("M58").Value = -1233.1842
$ws.Range("N58").Value = -23813160
$ws.Range("H107").Value = 1410.5
$ws.Range("I107").Value = 248.4375
$ws.Range("J107").Value = 3269.8
$ws.Range("K107").Value = 248.4375
$ws.Range("L107").Value = 3269.8
$ws.Range("M107").Value = 1671.5625
$ws.Range("N107").Value = -7109.8
$ws.Range("H136").Value = 8476651
$ws.Range("I136").Value = 1436.1842
$ws.Range("J136").Value = 23812754
$ws.Range("K136").Value = 4308.5526
$ws.Range("L136").Value = 71438262
$ws.Range("M136").Value = -1758.5526
$ws.Range("N136").Value = -71443362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 134274.14
$ws.Range("I4").Value = 222290
$ws.Range("J4").Value = 2250.3333
$ws.Range("K4").Value = 666870
$ws.Range("L4").Value = 6750.999899999999
$ws.Range("M4").Value = -666758
$ws.Range("N4").Value = -6974.999899999999
$ws.Range("H113").Value = 2321.9048
$ws.Range("I113").Value = 30003
$ws.Range("J113").Value = 937.85
$ws.Range("K113").Value = 90009
$ws.Range("L113").Value = 2813.55
$ws.Range("M113").Value = -87839
$ws.Range("N113").Value = -7153.55
$ws.Range("H131").Value = 982.6842
$ws.Range("I131").Value = 474.2857
$ws.Range("J131").Value = 1053.86
$ws.Range("K131").Value = 1422.8571
$ws.Range("L131").Value = 3161.58
$ws.Range("M131").Value = 3617.1429
$ws.Range("N131").Value = -13241.58

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3554.5454
$ws.Range("I113").Value = 3057.1428
$ws.Range("J113").Value = 4425
$ws.Range("K113").Value = 3057.1428
$ws.Range("L113").Value = 4425
$ws.Range("M113").Value = -887.1428000000001
$ws.Range("N113").Value = -8765
$ws.Range("H122").Value = 6773.4546
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 7350.8
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 22052.4
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -26952.4
$ws.Range("H126").Value = 3301.75
$ws.Range("I126").Value = 1813.3334
$ws.Range("K126").Value = 5440.0002
$ws.Range("M126").Value = -2970.0002
$ws.Range("H139").Value = 29590.908
$ws.Range("J139").Value = 29590.908
$ws.Range("L139").Value = 29590.908
$ws.Range("N139").Value = -39870.908

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4666.6665
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 4666.6665
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 4666.6665
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -4890.6665
$ws.Range("H100").Value = 2075.8823
$ws.Range("I100").Value = 1532.2222
$ws.Range("K100").Value = 1532.2222
$ws.Range("M100").Value = -991.2221999999999
$ws.Range("H122").Value = 3191.5386
$ws.Range("I122").Value = 2587.0588
$ws.Range("J122").Value = 4333.3335
$ws.Range("K122").Value = 7761.176399999999
$ws.Range("L122").Value = 13000.0005
$ws.Range("M122").Value = -5311.176399999999
$ws.Range("N122").Value = -17900.0005
$ws.Range("H126").Value = 4666.6665
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 4666.6665
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 13999.9995
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -18939.9995
$ws.Range("H132").Value = 2530.359
$ws.Range("I132").Value = 1411.44
$ws.Range("K132").Value = 4234.32
$ws.Range("M132").Value = -1704.32

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1070.4
$ws.Range("I100").Value = 1086.2858
$ws.Range("J100").Value = 1033.3334
$ws.Range("K100").Value = 2172.5716
$ws.Range("L100").Value = 2066.6668
$ws.Range("M100").Value = -1631.5716
$ws.Range("N100").Value = -3148.6668
